$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5984.301031857052
$ws.Range("B4").Value = 13191.30596193509
$ws.Range("B5").Value = 311335.9915490934
$ws.Range("B6").Value = 1393.455832459419
$ws.Range("B7").Value = 5984.301031857052
$ws.Range("B8").Value = 13191.30596193509
$ws.Range("B9").Value = 311335.9915490934
$ws.Range("B10").Value = 1393.455832459419
$ws.Range("B14").Value = 128.0105706272269
$ws.Range("B15").Value = 5864.615011217259
$ws.Range("B16").Value = 13059.39290231574
$ws.Range("B17").Value = 143.4281960832142
$ws.Range("B18").Value = 119.6860206370869
$ws.Range("B19").Value = 131.9130596193509
$ws.Range("B20").Value = 4.435923590202504
$ws.Range("B21").Value = 311335.9915490932
$ws.Range("B22").Value = 5984.301031854346
$ws.Range("B23").Value = 13191.30596193509
$ws.Range("B24").Value = 311335.9915490932
$ws.Range("B25").Value = 147.8641196734167
$ws.Range("B26").Value = 330511.5985428825
$ws.Range("B27").Value = 6653.885385303753
$ws.Range("B28").Value = 37347.81063534573
$ws.Range("B29").Value = 8216.518339776061
$ws.Range("B30").Value = 5805.968861105086
$ws.Range("B31").Value = 130.5939290231574
$ws.Range("B32").Value = 58.64615011217261
$ws.Range("B33").Value = 12928.79897329258
$ws.Range("B34").Value = 5864.615011217259
$ws.Range("B35").Value = 13059.39290231574
$ws.Range("B36").Value = 18924.00791353299
$ws.Range("B37").Value = 14571.48609342041
$ws.Range("B38").Value = 3205.726940552489
$ws.Range("B39").Value = 5805.968861105086
$ws.Range("B40").Value = 12928.79897329258
$ws.Range("B41").Value = 1393.455832459419
$ws.Range("B42").Value = 311592.0265529398
$ws.Range("B43").Value = 219.3279897594884
$ws.Range("B44").Value = 35244.12468267969
$ws.Range("B46").Value = 18078.77066563231
$ws.Range("B47").Value = 26938.09865539179
$ws.Range("B51").Value = 8306.0260272879
